$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("County")
$ws.Range("B1").Value = 160
[void]$ws.Range("B2").Select()
